$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2024-10-05"
$ws.Range("A6").ClearFormats()
$ws.Range("B6").Value = 0.02912
